# Bouton retour à l'étape de validation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style 4 = center + vertical-center + wrap (same style already used by the
#     existing header-row product cells). We reuse it for the whole product grid.

# Column A (A2:A16) - product / dish names
$ws.Range("A2").Value = "Petit choux de blé noir garni"
$ws.Range("A3").Value = "Blini garni"
$ws.Range("A4").Value = "Mini-roulées de blé noir"
$ws.Range("A5").Value = "Tuiles"
$ws.Range("A6").Value = "Triskels au chocolat"
$ws.Range("A7").Value = "Truffes"
$ws.Range("A8").Value = "Meringues"
$ws.Range("A9").Value = "Kig ha farz"
$ws.Range("A10").Value = "Potée de pouldrezic (aux choux)"
$ws.Range("A11").Value = "Potée Guérandaise (fèves, lard, saucisses)"
$ws.Range("A12").Value = "Frigousse de bœuf"
$ws.Range("A13").Value = "Cotriade ou matelote"
$ws.Range("A14").Value = "Poulet au cidre"
$ws.Range("A15").Value = "Jambon à l'os (environ 30 pers.)"
$ws.Range("A16").Value = "Buffet de crêpes : peut-être accompagné de garnitures (sucre, confitures, ...)"

# B2 and the existing menu-description cells C2:D4 get new / updated text
$ws.Range("B2").Value = "1 galette blé noir_2 crêpes froment"
$ws.Range("C2").Value = "Café, thé, jus de pommes ou raisins_Triskels au chocolat_Meringues"
$ws.Range("C3").Value = "Café, thé, jus de pommes ou raisins_Gâteau breton_Triskels au chocolat"
$ws.Range("C4").Value = "Café, thé, jus de pommes ou raisins_Gâteau breton_Triskels au chocolat_Crêpes roulées"
$ws.Range("D2").Value = "Cidre, vin blanc, jus de pommes_Assortiment de crêpes roulées salées"
$ws.Range("D3").Value = "Cidre, vin blanc, jus de raisins_Assortiment de lichouseries sucrées et salées"

# Remaining column B cells (B3:B10) - new cider / drinks menu rows
$ws.Range("B3").Value = "2 galette blé noir_2 crêpes froment"
$ws.Range("B4").Value = "galettes blé noir à volonté_crêpes froment à volonté"
$ws.Range("B5").Value = "1 galette blé noir_1 crêpes froment"
$ws.Range("B6").Value = "Assortiment de lichouseries_2 galettes blé noir_Salade_2 crêpes froment"
$ws.Range("B7").Value = "Cidre artisanal Kerné"
$ws.Range("B8").Value = "Cidre fermier Melenig"
$ws.Range("B9").Value = "Cidre Kerné"
$ws.Range("B10").Value = "Jus de pomme Kerné"

# Harmonise styling: rows 2-16 (A:D) all use the wrap+vcenter+center style (index 4
# in styles.xml), matching the original C/D header cells' style.
$ws.Range("A2:D16").WrapText = $true
$ws.Range("A2:D16").HorizontalAlignment = -4108
$ws.Range("A2:D16").VerticalAlignment = -4108

# Row heights: 43.2 for rows 2-4, 28.8 for row 6 and row 16 (matches autosize of the
# longer wrapped text), default elsewhere.
$ws.Rows.Item(2).RowHeight = 43.2
$ws.Rows.Item(3).RowHeight = 43.2
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 28.8
$ws.Rows.Item(16).RowHeight = 28.8

# Update the active selection to C11 (as recorded in the saved file)
$ws.Range("C11").Select()
